$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column A (Date) as text so date-like strings are stored as literal text,
# matching the inlineStr "YYYY-MM-DD" cells used throughout the rest of the sheet.
$ws.Range("A618:A638").NumberFormat = "@"

# Row 618
$ws.Range("A618").Value = "2024-08-28"
$ws.Range("C618").Value = 1939.099975585938
$ws.Range("D618").Value = 1719.449951171875
$ws.Range("E618").Value = 1078.800048828125
$ws.Range("F618").Value = 1847.050048828125
$ws.Range("G618").Value = 1737.900024414062
$ws.Range("H618").Value = 37879.2001953125
$ws.Range("I618").Value = 0
$ws.Range("J618").Value = 176.0689475795104

# Row 619
$ws.Range("A619").Value = "2024-08-29"
$ws.Range("C619").Value = 1933.349975585938
$ws.Range("D619").Value = 1751.849975585938
$ws.Range("E619").Value = 1061.300048828125
$ws.Range("F619").Value = 1857.849975585938
$ws.Range("G619").Value = 1726.550048828125
$ws.Range("H619").Value = 37876.40014648438
$ws.Range("I619").Value = -0.00007392048442647694
$ws.Range("J619").Value = 176.0559324776129

# Row 620
$ws.Range("A620").Value = "2024-08-30"
$ws.Range("C620").Value = 1943.699951171875
$ws.Range("D620").Value = 1753.25
$ws.Range("E620").Value = 1065.599975585938
$ws.Range("F620").Value = 1815.150024414062
$ws.Range("G620").Value = 1772.25
$ws.Range("H620").Value = 37943.19958496094
$ws.Range("I620").Value = 0.001763616347335551
$ws.Range("J620").Value = 176.3664275981758

# Row 621
$ws.Range("A621").Value = "2024-09-02"
$ws.Range("C621").Value = 1964.5
$ws.Range("D621").Value = 1806.650024414062
$ws.Range("E621").Value = 1050.949951171875
$ws.Range("F621").Value = 1766.300048828125
$ws.Range("G621").Value = 1749.5
$ws.Range("H621").Value = 38025.74987792969
$ws.Range("I621").Value = 0.002175628146063607
$ws.Range("J621").Value = 176.7501353620791

# Row 622
$ws.Range("A622").Value = "2024-09-03"
$ws.Range("C622").Value = 1941.25
$ws.Range("D622").Value = 1790.449951171875
$ws.Range("E622").Value = 1068.800048828125
$ws.Range("F622").Value = 1769.300048828125
$ws.Range("G622").Value = 1718.75
$ws.Range("H622").Value = 37895.55029296875
$ws.Range("I622").Value = -0.003423984678248408
$ws.Range("J622").Value = 176.144945606721

# Row 623
$ws.Range("A623").Value = "2024-09-04"
$ws.Range("C623").Value = 1922.449951171875
$ws.Range("D623").Value = 1785.25
$ws.Range("E623").Value = 1056.199951171875
$ws.Range("F623").Value = 1749.699951171875
$ws.Range("G623").Value = 1729.550048828125
$ws.Range("H623").Value = 37618.74926757812
$ws.Range("I623").Value = -0.007304314708473397
$ws.Range("J623").Value = 174.8583274897026

# Row 624
$ws.Range("A624").Value = "2024-09-05"
$ws.Range("C624").Value = 1933.150024414062
$ws.Range("D624").Value = 1790.550048828125
$ws.Range("E624").Value = 1074.900024414062
$ws.Range("F624").Value = 1722.900024414062
$ws.Range("G624").Value = 1720.75
$ws.Range("H624").Value = 37761.05065917969
$ws.Range("I624").Value = 0.003782725220059497
$ws.Range("J624").Value = 175.5197684950353

# Row 625
$ws.Range("A625").Value = "2024-09-06"
$ws.Range("C625").Value = 1901.849975585938
$ws.Range("D625").Value = 1756.099975585938
$ws.Range("E625").Value = 1112.650024414062
$ws.Range("F625").Value = 1730.300048828125
$ws.Range("G625").Value = 1715
$ws.Range("H625").Value = 37713.70007324219
$ws.Range("I625").Value = -0.001253953084220899
$ws.Range("J625").Value = 175.2996749399892

# Row 626
$ws.Range("A626").Value = "2024-09-09"
$ws.Range("C626").Value = 1894.650024414062
$ws.Range("D626").Value = 1746.75
$ws.Range("E626").Value = 1077.550048828125
$ws.Range("F626").Value = 1750.400024414062
$ws.Range("G626").Value = 1741.199951171875
$ws.Range("H626").Value = 37455.65051269531
$ws.Range("I626").Value = -0.006842329446480399
$ws.Range("J626").Value = 174.1002168121888

# Row 627
$ws.Range("A627").Value = "2024-09-10"
$ws.Range("C627").Value = 1912.300048828125
$ws.Range("D627").Value = 1779.099975585938
$ws.Range("E627").Value = 1091
$ws.Range("F627").Value = 1756.349975585938
$ws.Range("G627").Value = 1745.150024414062
$ws.Range("H627").Value = 37856.65014648438
$ws.Range("I627").Value = 0.01070598503296977
$ws.Range("J627").Value = 175.964131127617

# Row 628
$ws.Range("A628").Value = "2024-09-11"
$ws.Range("C628").Value = 1910.150024414062
$ws.Range("D628").Value = 1778.75
$ws.Range("E628").Value = 1077.849975585938
$ws.Range("F628").Value = 1789.349975585938
$ws.Range("G628").Value = 1782.650024414062
$ws.Range("H628").Value = 37910.79992675781
$ws.Range("I628").Value = 0.001430390170918655
$ws.Range("J628").Value = 176.2158284912161

# Row 629
$ws.Range("A629").Value = "2024-09-12"
$ws.Range("C629").Value = 1950.449951171875
$ws.Range("D629").Value = 1807.599975585938
$ws.Range("E629").Value = 1083.75
$ws.Range("F629").Value = 1838.050048828125
$ws.Range("G629").Value = 1812.75
$ws.Range("H629").Value = 38550.34973144531
$ws.Range("I629").Value = 0.01686985782212682
$ws.Range("J629").Value = 179.1885644638713

# Row 630
$ws.Range("A630").Value = "2024-09-13"
$ws.Range("C630").Value = 1944.099975585938
$ws.Range("D630").Value = 1812.800048828125
$ws.Range("E630").Value = 1089.699951171875
$ws.Range("F630").Value = 1826.050048828125
$ws.Range("G630").Value = 1814.099975585938
$ws.Range("H630").Value = 38552.5498046875
$ws.Range("I630").Value = 0.00005707012407186833
$ws.Range("J630").Value = 179.1987907774775

# Row 631
$ws.Range("A631").Value = "2024-09-16"
$ws.Range("C631").Value = 1950.25
$ws.Range("D631").Value = 1811.849975585938
$ws.Range("E631").Value = 1094.650024414062
$ws.Range("F631").Value = 1757.849975585938
$ws.Range("G631").Value = 1797.199951171875
$ws.Range("H631").Value = 38385.89990234375
$ws.Range("I631").Value = -0.004322668751821118
$ws.Range("J631").Value = 178.4241737642195

# Row 632
$ws.Range("A632").Value = "2024-09-17"
$ws.Range("C632").Value = 1952.550048828125
$ws.Range("D632").Value = 1813.75
$ws.Range("E632").Value = 1080.300048828125
$ws.Range("F632").Value = 1741.150024414062
$ws.Range("G632").Value = 1848.5
$ws.Range("H632").Value = 38346.90075683594
$ws.Range("I632").Value = -0.001015975803798501
$ws.Range("J632").Value = 178.2428991208623

# Row 633
$ws.Range("A633").Value = "2024-09-18"
$ws.Range("C633").Value = 1892.150024414062
$ws.Range("D633").Value = 1756.5
$ws.Range("E633").Value = 1065.800048828125
$ws.Range("F633").Value = 1727.25
$ws.Range("G633").Value = 1805.599975585938
$ws.Range("H633").Value = 37454.75048828125
$ws.Range("I633").Value = -0.023265250931541
$ws.Range("J633").Value = 174.0960333460501

# Row 634
$ws.Range("A634").Value = "2024-09-19"
$ws.Range("C634").Value = 1894.199951171875
$ws.Range("D634").Value = 1736.5
$ws.Range("E634").Value = 1060.75
$ws.Range("F634").Value = 1676.449951171875
$ws.Range("G634").Value = 1877.449951171875
$ws.Range("H634").Value = 37317.94946289062
$ws.Range("I634").Value = -0.003652434567236724
$ws.Range("J634").Value = 173.4601589758382

# Row 635
$ws.Range("A635").Value = "2024-09-20"
$ws.Range("C635").Value = 1905.75
$ws.Range("D635").Value = 1760.050048828125
$ws.Range("E635").Value = 1114.699951171875
$ws.Range("F635").Value = 1662
$ws.Range("G635").Value = 1931.449951171875
$ws.Range("H635").Value = 38001.24975585938
$ws.Range("I635").Value = 0.01831023147850691
$ws.Range("J635").Value = 176.6362546389844

# Row 636
$ws.Range("A636").Value = "2024-09-23"
$ws.Range("C636").Value = 1896.449951171875
$ws.Range("D636").Value = 1752.800048828125
$ws.Range("E636").Value = 1106.699951171875
$ws.Range("F636").Value = 1692.900024414062
$ws.Range("G636").Value = 1920.400024414062
$ws.Range("H636").Value = 37915.79968261719
$ws.Range("I636").Value = -0.002248612184892999
$ws.Range("J636").Value = 176.2390682045093

# Row 637
$ws.Range("A637").Value = "2024-09-24"
$ws.Range("C637").Value = 1898.599975585938
$ws.Range("D637").Value = 1775.599975585938
$ws.Range("E637").Value = 1098.5
$ws.Range("F637").Value = 1660.900024414062
$ws.Range("G637").Value = 1838.75
$ws.Range("H637").Value = 37717.7998046875
$ws.Range("I637").Value = -0.00522209420840627
$ws.Range("J637").Value = 175.3187311871436

# Row 638
$ws.Range("A638").Value = "2024-09-25"
$ws.Range("C638").Value = 1895.300048828125
$ws.Range("D638").Value = 1782.400024414062
$ws.Range("E638").Value = 1088.599975585938
$ws.Range("F638").Value = 1654.75
$ws.Range("G638").Value = 1722.050048828125
$ws.Range("H638").Value = 37400.95031738281
$ws.Range("I638").Value = -0.008400529430280024
$ws.Range("J638").Value = 173.8459610261267
